# Commit: "add conversion to arctos"
#
# Upstream, 8 new rows were inserted earlier in the master Lewis dataset,
# which shifts the running TPTID counter for every taxon that was
# reviewed after that point. This review sheet (Lewis_review) stamps each
# row with the TPTID at the time of review (column I, "TPTID") and the
# matching taxonID string (column J, "Lewis<TPTID>"), so every one of
# those stamped IDs needs to be renumbered down by 8 to stay in sync.
#
# Concretely: for every data row, TPTID := TPTID - 8, and taxonID :=
# "Lewis" + TPTID.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NB: plain `.Value` on this host reflects back a property signature
# string instead of invoking the getter/setter for Range; `.Value2`
# behaves correctly for both read and write, so use that throughout.

$lastRow = $ws.Cells.Item($ws.Rows.Count, 9).End(-4162).Row   # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $tptid = $ws.Cells.Item($r, 9).Value2
    if ($tptid -eq $null) { continue }

    $newId = [int]$tptid - 8

    $ws.Cells.Item($r, 9).Value2  = $newId
    $ws.Cells.Item($r, 10).Value2 = "Lewis" + $newId
}
